$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -4
    3  = 3
    4  = 3
    5  = 1
    6  = 5
    7  = 3
    8  = -2
    9  = 9
    10 = -2
    12 = 2
    13 = 4
    14 = 2
    15 = -4
    17 = 3
    18 = -2
    19 = -1
    20 = 1
    21 = -1
    22 = -5
    23 = -4
    24 = 2
    25 = 1
    26 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
